$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Programs" SQL query in cell B2 -------------------------
# The query now derives "Website" from a CASE expression over
# prg.program_link / prg.program_acronym instead of a plain prg.website
# column reference.
$newProgramsQuery = "SELECT DISTINCT `r`n    prg.program_name AS ""Program"",`r`n  CASE`r`n    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym`r`n        ELSE prg.program_link`r`n    END  AS ""Website"",`r`n    prg.focus_area AS ""Focus Area"",`r`n    prg.cancer_type AS ""Cancer Type"",`r`n CASE `r`n        WHEN prg.data_link IS NOT NULL THEN prg.website       `r`n        ELSE prg.data_link`r`n    END AS ""Data Location Details""`r`nFROM `r`n    df_program prg`r`nWHERE `r`n     prg.cancer_type LIKE '%Leukemia%'`r`nORDER BY `r`n   lower(prg.program_name) ASC`r`nLIMIT 100;"

# Touching Font.ThemeColor (re-asserting the same theme color already in
# effect) mirrors what happened in the authored workbook: it causes the
# cell's effective style to be re-materialized as a fresh (but visually
# identical - size 12, wrapped) style entry rather than reusing the old
# one, matching the style-index bump seen on B2 in the saved file.
$ws.Range("B2").Font.ThemeColor = 1
$ws.Range("B2").Value = $newProgramsQuery

# --- Update the view state left behind by the edit -----------------------
# The author ended the session scrolled down with C8 selected.
$ws.Range("C8").Select() | Out-Null
